# Infrastructure LLD.xlsx - "Added Redis and blob storage"
#
# Re-lays out Sheet1: three new Q&A rows are inserted near the top (new
# row 3, pushing the existing SANDBOX-endpoint question down to row 5
# and everything below it down by three), a third column (C) is added
# holding short owner/notes tags next to several existing and new
# questions, and three new rows are added near the bottom of the
# tracker block (ALI / Devop engineer, Jevgenijs, and Next steps /
# Ali to catch up and start), plus three blank trailer rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlTop = -4160
$xlLeft = -4131

# ---------------------------------------------------------------------
# New column C
# ---------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 46.7109375

# ---------------------------------------------------------------------
# Rows 4-6 previously held A-column placeholder cells (style only, no
# value) belonging to rows that are being pushed further down the
# sheet. The new rows 4-6 have no A-column cell at all, so drop that
# stale formatting before writing the new B/C content.
# ---------------------------------------------------------------------
$ws.Range("A4").Clear()
$ws.Range("A5").Clear()
$ws.Range("A6").Clear()

# Row 3 (new): "Review available templates" question
$ws.Range("B3").Value = "Review available templates:`nhttps://dev.azure.com/defragovuk/DEFRA-DEVOPS-COMMON/_git/Defra.Infrastructure.Common?path=/templates/Microsoft.Network"

# Row 4 (new): "Git hub actions..." question + owner note
$ws.Range("B4").Value = "Git hub actions or dev ops pipelines?"
$ws.Range("C4").Value = "Pieplines and current project"

# Row 5: existing SANDBOX service-endpoint question, now here, + note
$ws.Range("B5").Value = "Is it possible to get a service endpoint created in the SANDBOX to create the end to end devops pipeline"
$ws.Range("C5").Value = "Ali to do "

# Row 6 (new): Role assignments / rbac question + owner note
$ws.Range("B6").Value = "Role assignments - how do we handle? I intend to use the rbac model.`nIdeally want to set all of these permissions in the pipeline, e.g webapp can get a secret from the keyvautl"
$ws.Range("C6").Value = "App registration and rbac"

# Row 7: virtual network question (shifted down from old row 4)
$ws.Range("B7").Value = "A virtual network with:`n- A private endpoint subnet`n- An app service subnet with:`ndelegation to 'Microsoft.Web\serverFarms'`nService endpoint enabled for 'Microsoft.Web'"

# Row 8: resource ID's question (shifted down from old row 5)
$ws.Range("B8").Value = "Need to know resource ID's of virtual networks and subnets"

# Row 9: Azure policy / log analytics question + new owner note
$ws.Range("B9").Value = "Will Azure policy enforce logging to centralised log analytics instance?"
$ws.Range("C9").Value = "Resources Will be connected to central LA.`n"
$ws.Range("C9").WrapText = $true

# Row 10: resource ID of centralised log analytics instance
$ws.Range("B10").Value = "We will require the resource ID of the centralised log analytics instance for each environment"

# Row 11: Finalise naming convention
$ws.Range("B11").Value = "Finalise naming convention"

# Row 12: private dns question
$ws.Range("B12").Value = "Need to know how private dns will be handled. `nIs this done outside of the project pipeline"

# Row 13: service connection question
$ws.Range("B13").Value = "Will need a service connection created in the project devops instance for SND, DEV and TEST"

# Row 14: DEV OPS agent question
$ws.Range("B14").Value = "Will need a DEV OPS agent that can access private endpoints for DEV and TEST"

# Row 15: mechanism to test access + new owner note
$ws.Range("B15").Value = "We will need a mechanism to test access to web apps"
$ws.Range("C15").Value = "Front end ingress point"

# Row 16: egress connectivity question + new owner note
$ws.Range("B16").Value = "How do we manage egress connectivity to other azure resources, internet addresses etc"
$ws.Range("C16").Value = "All egress goes via palo alto."

# Row 17: API Manager question + new owner note
$ws.Range("B17").Value = "How do we manage the API Manager configuration?"
$ws.Range("C17").Value = "TBD. CCOE engineer allocated to work through."

# Row 18: Azure Front door question + new owner note
$ws.Range("B18").Value = "How do we manage Azure Front door?`n- Lower environments`n- Private access to exposed services (assume IP restrictions in WAF policies)`n- How can we view firewall logs?"
$ws.Range("C18").Value = "F5 silverline"

# Row 19: Entra ID app registrations question + new owner note
$ws.Range("B19").Value = "How do we manage the Entra ID app registrations?"
$ws.Range("C19").Value = "TBD. CCOE engineer allocated to work through."

# ---------------------------------------------------------------------
# Row 23: new - ALI / Devop engineer
# ---------------------------------------------------------------------
$ws.Range("A23").Value = "ALI "
$ws.Range("B23").Value = "Devop engineer"

# Row 24: new - Jevgenijs
$ws.Range("A24").Value = "Jevgenijs"

# Row 26: new - Next steps / Ali to catch up and start
$ws.Range("A26").Value = "Next steps"
$ws.Range("B26").Value = "Ali to catch up and start "

# ---------------------------------------------------------------------
# Rows 32-34 switch from the vertical-top-only placeholder style used
# further down the sheet to the left+top style used by the rows above
# them (they already have vertical=top from style 2; only horizontal
# alignment needs to change to pick up the existing style 3 xf).
# ---------------------------------------------------------------------
$ws.Range("A32").HorizontalAlignment = $xlLeft
$ws.Range("A33").HorizontalAlignment = $xlLeft
$ws.Range("A34").HorizontalAlignment = $xlLeft

# ---------------------------------------------------------------------
# Rows 38-40: new trailer rows, matching the vertical-top style used
# by rows 35-37 immediately above them.
# ---------------------------------------------------------------------
$ws.Range("A38").VerticalAlignment = $xlTop
$ws.Range("A39").VerticalAlignment = $xlTop
$ws.Range("A40").VerticalAlignment = $xlTop

# ---------------------------------------------------------------------
# Row heights: wrapped multi-line content needs an explicit height
# taller than the 15pt default; a couple of rows whose old content
# had an explicit height now hold single-line content and drop back
# to the default (AutoFit clears the explicit height cleanly).
# ---------------------------------------------------------------------
$ws.Rows(3).RowHeight = 45
$ws.Rows(4).AutoFit()
$ws.Rows(5).RowHeight = 30
$ws.Rows(6).RowHeight = 45
$ws.Rows(7).RowHeight = 75
$ws.Rows(12).RowHeight = 30
$ws.Rows(13).RowHeight = 30
$ws.Rows(15).AutoFit()
$ws.Rows(18).RowHeight = 60

# ---------------------------------------------------------------------
# Final selection, matching the saved workbook view.
# ---------------------------------------------------------------------
$ws.Range("B26").Select()
